$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is numeric-looking need to be forced to Text
# format first, otherwise Excel would auto-convert them to numbers (losing
# e.g. trailing zeros) instead of keeping them as plain text like the source feed.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated Price (D) and Volume(1h) (E) values per row.
$ws.Range("D2").Value = "29.099.05"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.835.61"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "243.28"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "0.6173"
$ws.Range("E6").Value = "  -2.03%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").Value = "0.07460"
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").Value = "0.2920"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("D11").Value = "0.07700"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "1.839.39"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").Value = "5.003"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "82.60"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").Value = "0.000009268"
$ws.Range("E16").Value = "  -3.88%  "
$ws.Range("D17").Value = "5.927"
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("D18").Value = "29.070.35"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "2.088.51"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "230.83"
$ws.Range("E20").Value = "  +2.01%  "
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").Value = "7.175"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("D25").Value = "160.33"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").Value = "0.1386"
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("D27").Value = "8.508"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "17.79"
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D29").Value = "1.499"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("D31").Value = "4.134"
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("D32").Value = "0.05511"
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("D33").Value = "1.212"
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("D34").Value = "0.7470"
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("D35").Value = "1.837"
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("D36").Value = "1.141"
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("D38").Value = "2.770"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("D39").Value = "1.217.15"
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").Value = "6.470"
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("D42").Value = "0.8960"
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").Value = "101.81"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").Value = "1.988.03"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "65.56"
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("D47").Value = "0.00000000123"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").Value = "0.5093"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").Value = "0.4065"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "9.110"
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("D51").Value = "0.07261"
$ws.Range("E51").Value = "  +10.31%  "
